$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('ECs', 'Mfge8', 'Pdgfrb', 'ECs', 3, 1, 15.12165933333333, 45.364978, 0.1696222886509932, 0.1696222886509932, 2, 0.6666666666666666, 4.389289333333333, 13.167868, 0.03995844413671427, 0.03995844413671427, 66.37333801410044, 597.360042126904, 0.006777842745402337, 0.006777842745402336),
    @('ECs', 'Mfge8', 'Pdgfrb', 'FAPs', 3, 1, 15.12165933333333, 45.364978, 0.1696222886509932, 0.1696222886509932, 3, 1, 57.89762366666667, 173.692871, 0.5270782546422108, 0.5270782546422108, 875.5081412968709, 7879.573271671838, 0.0894042198505828, 0.08940421985058279),
    @('ECs', 'Mfge8', 'Pdgfrb', 'M2', 3, 1, 15.12165933333333, 45.364978, 0.1696222886509932, 0.1696222886509932, 3, 1, 0.129282, 0.387846, 0.001176934848120294, 0.001176934848120294, 1.954958361932, 17.594625257388, 0.0001996343825312735, 0.0001996343825312734),
    @('ECs', 'Mfge8', 'Pdgfrb', 'sCs', 3, 1, 15.12165933333333, 45.364978, 0.1696222886509932, 0.1696222886509932, 3, 1, 47.43015733333333, 142.290472, 0.4317863663729547, 0.4317863663729548, 717.2226813210684, 6455.004131889616, 0.07324059167247685, 0.07324059167247685),
    @('FAPs', 'Mfge8', 'Pdgfrb', 'ECs', 3, 1, 30.93224, 92.79672, 0.3469723279862584, 0.3469723279862583, 2, 0.6666666666666666, 4.389289333333333, 13.167868, 0.03995844413671427, 0.03995844413671427, 135.7705510881067, 1221.93495979296, 0.01386447438482461, 0.01386447438482461),
    @('FAPs', 'Mfge8', 'Pdgfrb', 'FAPs', 3, 1, 30.93224, 92.79672, 0.3469723279862584, 0.3469723279862583, 3, 1, 57.89762366666667, 173.692871, 0.5270782546422108, 0.5270782546422108, 1790.903190687013, 16118.12871618312, 0.1828815690441418, 0.1828815690441417),
    @('FAPs', 'Mfge8', 'Pdgfrb', 'M2', 3, 1, 30.93224, 92.79672, 0.3469723279862584, 0.3469723279862583, 3, 1, 0.129282, 0.387846, 0.001176934848120294, 0.001176934848120294, 3.99898185168, 35.99083666512, 0.000408363824140452, 0.0004083638241404519),
    @('FAPs', 'Mfge8', 'Pdgfrb', 'sCs', 3, 1, 30.93224, 92.79672, 0.3469723279862584, 0.3469723279862583, 3, 1, 47.43015733333333, 142.290472, 0.4317863663729547, 0.4317863663729548, 1467.121009872427, 13204.08908885184, 0.1498179207331516, 0.1498179207331516),
    @('M2', 'Mfge8', 'Pdgfrb', 'ECs', 3, 1, 4.447555666666666, 13.342667, 0.04988900718188559, 0.04988900718188558, 2, 0.6666666666666666, 4.389289333333333, 13.167868, 0.03995844413671427, 0.03995844413671427, 19.52160864710622, 175.694477823956, 0.001993487106513512, 0.001993487106513512),
    @('M2', 'Mfge8', 'Pdgfrb', 'FAPs', 3, 1, 4.447555666666666, 13.342667, 0.04988900718188559, 0.04988900718188558, 3, 1, 57.89762366666667, 173.692871, 0.5270782546422108, 0.5270782546422108, 257.5029042252174, 2317.526138026957, 0.02629541083126097, 0.02629541083126097),
    @('M2', 'Mfge8', 'Pdgfrb', 'M2', 3, 1, 4.447555666666666, 13.342667, 0.04988900718188559, 0.04988900718188558, 3, 1, 0.129282, 0.387846, 0.001176934848120294, 0.001176934848120294, 0.5749888916979999, 5.174900025282, 0.00005871611109048479, 0.00005871611109048478),
    @('M2', 'Mfge8', 'Pdgfrb', 'sCs', 3, 1, 4.447555666666666, 13.342667, 0.04988900718188559, 0.04988900718188558, 3, 1, 47.43015733333333, 142.290472, 0.4317863663729547, 0.4317863663729548, 210.9482650187582, 1898.534385168824, 0.02154139313302062, 0.02154139313302062),
    @('sCs', 'Mfge8', 'Pdgfrb', 'ECs', 3, 1, 38.64755633333333, 115.942669, 0.4335163761808628, 0.4335163761808628, 2, 0.6666666666666666, 4.389289333333333, 13.167868, 0.03995844413671427, 0.03995844413671427, 169.6353067732991, 1526.717760959692, 0.01732263989997382, 0.01732263989997382),
    @('sCs', 'Mfge8', 'Pdgfrb', 'FAPs', 3, 1, 38.64755633333333, 115.942669, 0.4335163761808628, 0.4335163761808628, 3, 1, 57.89762366666667, 173.692871, 0.5270782546422108, 0.5270782546422108, 2237.601672223633, 20138.4150500127, 0.2284970549162253, 0.2284970549162252),
    @('sCs', 'Mfge8', 'Pdgfrb', 'M2', 3, 1, 38.64755633333333, 115.942669, 0.4335163761808628, 0.4335163761808628, 3, 1, 0.129282, 0.387846, 0.001176934848120294, 0.001176934848120294, 4.996433377886, 44.96790040097401, 0.0005102205303580842, 0.0005102205303580841),
    @('sCs', 'Mfge8', 'Pdgfrb', 'sCs', 3, 1, 38.64755633333333, 115.942669, 0.4335163761808628, 0.4335163761808628, 3, 1, 47.43015733333333, 142.290472, 0.4317863663729547, 0.4317863663729548, 1833.059677438863, 16497.53709694977, 0.1871864608343057, 0.1871864608343057)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $row[$j]
    }
}
